# Updates cryptos list cell values to reflect latest price/volume data,
# matching the upstream commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.556.80"
$ws.Range("E2").Value = "'  +1.17%  "

$ws.Range("D3").Value = "'1.652.51"
$ws.Range("E3").Value = "'  +2.22%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.30%  "

$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "'  +0.26%  "

$ws.Range("D6").Value = "'301.07"
$ws.Range("E6").Value = "'  -0.40%  "

$ws.Range("D7").Value = "'0.3789"
$ws.Range("E7").Value = "'  +1.35%  "

$ws.Range("D8").Value = "'50.80"
$ws.Range("E8").Value = "'  -2.04%  "

$ws.Range("D9").Value = "'0.3556"
$ws.Range("E9").Value = "'  +0.44%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.227"
$ws.Range("E10").Value = "'  +1.08%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.08105"
$ws.Range("E11").Value = "'  -0.35%  "

$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "'  +0.32%  "

$ws.Range("D13").Value = "'22.13"
$ws.Range("E13").Value = "'  +0.22%  "

$ws.Range("D14").Value = "'6.421"
$ws.Range("E14").Value = "'  +0.01%  "

$ws.Range("D15").Value = "'7.447"
$ws.Range("E15").Value = "'  +2.87%  "

$ws.Range("D16").Value = "'0.00001200"
$ws.Range("E16").Value = "'  -1.31%  "

$ws.Range("D17").Value = "'1.646.40"
$ws.Range("E17").Value = "'  +1.95%  "

$ws.Range("D18").Value = "'97.35"
$ws.Range("E18").Value = "'  +2.63%  "

$ws.Range("D19").Value = "'0.06991"
$ws.Range("E19").Value = "'  +1.01%  "

$ws.Range("D20").Value = "'6.763"
$ws.Range("E20").Value = "'  +2.68%  "

$ws.Range("D21").Value = "'17.52"
$ws.Range("E21").Value = "'  +1.13%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.15%  "

$ws.Range("D23").Value = "'12.66"
$ws.Range("E23").Value = "'  +2.09%  "

$ws.Range("D24").Value = "'23.580.49"
$ws.Range("E24").Value = "'  +1.31%  "

$ws.Range("D25").Value = "'2.479"
$ws.Range("E25").Value = "'  -1.42%  "

$ws.Range("D26").Value = "'2.945"
$ws.Range("E26").Value = "'  -3.60%  "

$ws.Range("D27").Value = "'21.10"
$ws.Range("E27").Value = "'  +0.74%  "

$ws.Range("D28").Value = "'152.26"
$ws.Range("E28").Value = "'  -0.17%  "

$ws.Range("D29").Value = "'5.217"
$ws.Range("E29").Value = "'  +0.82%  "

$ws.Range("D30").Value = "'133.13"
$ws.Range("E30").Value = "'  +0.13%  "

$ws.Range("D31").Value = "'1.843.53"
$ws.Range("E31").Value = "'  +2.64%  "

$ws.Range("D32").Value = "'6.998"
$ws.Range("E32").Value = "'  +7.63%  "

$ws.Range("D33").Value = "'2.153"
$ws.Range("E33").Value = "'  +6.49%  "

$ws.Range("D34").Value = "'11.86"
$ws.Range("E34").Value = "'  +1.18%  "

$ws.Range("D35").Value = "'1.039"
$ws.Range("E35").Value = "'  -4.93%  "

$ws.Range("D36").Value = "'0.02749"
$ws.Range("E36").Value = "'  +0.69%  "

$ws.Range("D37").Value = "'0.08720"
$ws.Range("E37").Value = "'  -0.14%  "

$ws.Range("D38").Value = "'6.020"
$ws.Range("E38").Value = "'  +2.33%  "

$ws.Range("D39").Value = "'0.2460"
$ws.Range("E39").Value = "'  -0.09%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06898"
$ws.Range("E40").Value = "'  -0.42%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'13.12"
$ws.Range("E41").Value = "'  +4.66%  "

$ws.Range("D42").Value = "'0.6954"
$ws.Range("E42").Value = "'  +0.98%  "

$ws.Range("D43").Value = "'1.321"
$ws.Range("E43").Value = "'  -0.02%  "

$ws.Range("D44").Value = "'15.81"
$ws.Range("E44").Value = "'  +2.23%  "

$ws.Range("D45").Value = "'0.6455"
$ws.Range("E45").Value = "'  +1.91%  "

$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "'  +0.23%  "

$ws.Range("D47").Value = "'2.275"
$ws.Range("E47").Value = "'  +0.74%  "

$ws.Range("D48").Value = "'3.933"
$ws.Range("E48").Value = "'  -0.13%  "

$ws.Range("D49").Value = "'0.07873"
$ws.Range("E49").Value = "'  -0.03%  "

$ws.Range("D50").Value = "'126.93"
$ws.Range("E50").Value = "'  +0.07%  "

$ws.Range("D51").Value = "'1.184"
$ws.Range("E51").Value = "'  +1.84%  "
